#
# Adds the newly-processed paper "Early intervention in psychosis in
# Portugal: Where are we?" (Coentre R. Levy P., 2020) as entry #28 to
# the tracker workbook.
#
# Sheet "Summary" gets a brand-new data row (row 30) with the full
# record. The other five sheets ("Charactheristics", "Techniques",
# "Metrics", "Problems", "Citations") only get the index/title stub
# for the new paper (a row with just the paper number + title), which
# pushes their trailing "Sum Value"/"Count Times" aggregate rows down
# by one row.

$wb = $excel.ActiveWorkbook

$paperTitle  = "Early intervention in psychosis in Portugal: Where are we?"
$paperAuthor = "Coentre R. Levy P."
$paperDate   = "2020"
$dateSerial  = 44514

# ---------------------------------------------------------------------
# Sheet 1: Summary - append the full new record as row 30
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Summary")

# Copy formatting from the row above for the numeric index + date cells
# so no new cell styles get fabricated.
$ws1.Range("A29").Copy($ws1.Range("A30"))
$ws1.Range("A30").Value = 28

$ws1.Range("B30").Value = $paperTitle
$ws1.Range("C30").Value = $paperAuthor

# Column D holds the paper's publication year as text (shared string),
# not a number - reuse an existing cell already containing "2020" so it
# keeps the shared-string type instead of being coerced to a numeric
# value.
$ws1.Range("D15").Copy($ws1.Range("D30"))

$ws1.Range("E29").Copy($ws1.Range("E30"))
$ws1.Range("E30").Value = $dateSerial
$ws1.Range("F29").Copy($ws1.Range("F30"))
$ws1.Range("F30").Value = $dateSerial
$ws1.Range("G29").Copy($ws1.Range("G30"))
$ws1.Range("G30").Value = 0

# ---------------------------------------------------------------------
# Sheets 2-6: insert a stub row (index + title only) right before the
# trailing "Sum Value" / "Count Times" rows, shifting those down by one.
# ---------------------------------------------------------------------
$stubSheets = @("Charactheristics", "Techniques", "Metrics", "Problems", "Citations")

foreach ($name in $stubSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(30).Insert()
    $ws.Range("A29").Copy($ws.Range("A30"))
    $ws.Range("A30").Value = 28
    $ws.Range("B30").Value = $paperTitle
}
